$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1054.1177
$ws.Range("I80").Value = 778.5
$ws.Range("J80").Value = 1204.4546
$ws.Range("K80").Value = 2335.5
$ws.Range("L80").Value = 3613.3638
$ws.Range("M80").Value = -1337.5
$ws.Range("N80").Value = -5609.3638

$ws.Range("H83").Value = 1054.1177
$ws.Range("I83").Value = 778.5
$ws.Range("J83").Value = 1204.4546
$ws.Range("K83").Value = 7006.5
$ws.Range("L83").Value = 10840.0914
$ws.Range("M83").Value = -2014.5
$ws.Range("N83").Value = -20824.0914

$ws.Range("H98").Value = 802.38464
$ws.Range("I98").Value = 802.38464
$ws.Range("K98").Value = 802.38464
$ws.Range("M98").Value = 695.61536

$ws.Range("H122").Value = 802.38464
$ws.Range("I122").Value = 802.38464
$ws.Range("K122").Value = 2407.15392
$ws.Range("M122").Value = 42.84608000000026

$ws.Range("H137").Value = 3382.75
$ws.Range("I137").Value = 2076.2727
$ws.Range("K137").Value = 6228.8181
$ws.Range("M137").Value = -3678.8181

$ws.Range("H140").Value = 35000
$ws.Range("I140").Value = 35000
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 35000
$ws.Range("L140").ClearContents()
$ws.Range("N140").Value = 0
$ws.Range("M140").Value = -29820

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 1124.6666
$ws.Range("I19").Value = 1537.5
$ws.Range("K19").Value = 1537.5
$ws.Range("M19").Value = -1308.5

$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").ClearContents()
$ws.Range("N103").Value = 0

$ws.Range("H122").Value = 1897.8
$ws.Range("I122").Value = 1897.8
$ws.Range("K122").Value = 5693.4
$ws.Range("M122").Value = -3243.4

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").ClearContents()
$ws.Range("N128").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 10664
$ws.Range("I20").Value = 11098.667
$ws.Range("K20").Value = 11098.667
$ws.Range("M20").Value = -10851.667

$ws.Range("H92").Value = 10000
$ws.Range("J92").Value = 10000
$ws.Range("L92").Value = 10000
$ws.Range("N92").Value = -14992

$ws.Range("H105").Value = 1605.5
$ws.Range("I105").Value = 1307.5
$ws.Range("K105").Value = 1307.5
$ws.Range("M105").Value = 439.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4480.1304
$ws.Range("I31").Value = 1880.9678
$ws.Range("K31").Value = 1880.9678
$ws.Range("M31").Value = -1585.9678

$ws.Range("H34").Value = 4480.1304
$ws.Range("I34").Value = 1880.9678
$ws.Range("K34").Value = 1880.9678
$ws.Range("M34").Value = -1678.9678

$ws.Range("H99").Value = 3447.9312
$ws.Range("I99").Value = 3423.3809
$ws.Range("J99").Value = 3512.375
$ws.Range("K99").Value = 3423.3809
$ws.Range("L99").Value = 3512.375
$ws.Range("M99").Value = -1925.3809
$ws.Range("N99").Value = -6508.375

$ws.Range("H107").Value = 720.6667
$ws.Range("I107").Value = 736
$ws.Range("J107").Value = 644
$ws.Range("K107").Value = 736
$ws.Range("L107").Value = 644
$ws.Range("M107").Value = 1184
$ws.Range("N107").Value = -4484

$ws.Range("H126").Value = 3447.9312
$ws.Range("I126").Value = 3423.3809
$ws.Range("J126").Value = 3512.375
$ws.Range("K126").Value = 10270.1427
$ws.Range("L126").Value = 10537.125
$ws.Range("M126").Value = -7800.1427
$ws.Range("N126").Value = -15477.125

$ws.Range("H130").Value = 50000
$ws.Range("J130").Value = 50000
$ws.Range("L130").Value = 50000
$ws.Range("N130").Value = -60040

$ws.Range("H132").Value = 3597.7273
$ws.Range("I132").Value = 2958
$ws.Range("K132").Value = 8874
$ws.Range("M132").Value = -6344

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 39
$ws.Range("J2").Value = 45.18182
$ws.Range("L2").Value = 271.09092
$ws.Range("N2").Value = -497.09092

$ws.Range("H25").Value = 95
$ws.Range("I25").Value = 95
$ws.Range("K25").Value = 285
$ws.Range("M25").Value = -116

$ws.Range("H30").Value = 95
$ws.Range("I30").Value = 95
$ws.Range("K30").Value = 285
$ws.Range("M30").Value = -183

$ws.Range("H58").Value = 2501.25
$ws.Range("I58").Value = 2005
$ws.Range("J58").Value = 2666.6667
$ws.Range("K58").Value = 6015
$ws.Range("L58").Value = 8000.000100000001
$ws.Range("M58").Value = -5887
$ws.Range("N58").Value = -8256.000100000001

$ws.Range("H116").Value = 1965.4
$ws.Range("I116").Value = 2082
$ws.Range("J116").Value = 1499
$ws.Range("K116").Value = 6246
$ws.Range("L116").Value = 4497
$ws.Range("M116").Value = -2804
$ws.Range("N116").Value = -11381

$ws.Range("H132").Value = 2992.4375
$ws.Range("I132").Value = 2101
$ws.Range("K132").Value = 18909
$ws.Range("M132").Value = -16379

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").ClearContents()
$ws.Range("N108").Value = 0

$ws.Range("H113").Value = 7549.643
$ws.Range("I113").Value = 5527.857
$ws.Range("K113").Value = 5527.857
$ws.Range("M113").Value = -3357.857

$ws.Range("H122").Value = 503256.7
$ws.Range("I122").Value = 716114.4399999999
$ws.Range("J122").Value = 6588.6665
$ws.Range("K122").Value = 2148343.32
$ws.Range("L122").Value = 19765.9995
$ws.Range("M122").Value = -2145893.32
$ws.Range("N122").Value = -24665.9995

$ws.Range("H126").Value = 5012
$ws.Range("I126").Value = 5012
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 15036
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -12566

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7399.5
$ws.Range("I7").Value = 6519.4
$ws.Range("J7").Value = 8866.333000000001
$ws.Range("K7").Value = 6519.4
$ws.Range("L7").Value = 8866.333000000001
$ws.Range("M7").Value = -6407.4
$ws.Range("N7").Value = -9090.333000000001

$ws.Range("H46").Value = 8698.777
$ws.Range("J46").Value = 6077
$ws.Range("L46").Value = 6077
$ws.Range("N46").Value = -6453

$ws.Range("H82").Value = 6236.364
$ws.Range("I82").Value = 4930.2
$ws.Range("K82").Value = 4930.2
$ws.Range("M82").Value = -4569.2

$ws.Range("H85").Value = 6236.364
$ws.Range("I85").Value = 4930.2
$ws.Range("K85").Value = 4930.2
$ws.Range("M85").Value = -3682.2

$ws.Range("H126").Value = 7399.5
$ws.Range("I126").Value = 6519.4
$ws.Range("J126").Value = 8866.333000000001
$ws.Range("K126").Value = 19558.2
$ws.Range("L126").Value = 26598.999
$ws.Range("M126").Value = -17088.2
$ws.Range("N126").Value = -31538.999

$ws.Range("H132").Value = 3517.5667
$ws.Range("I132").Value = 2531.182
$ws.Range("K132").Value = 7593.545999999999
$ws.Range("M132").Value = -5063.545999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 30000
$ws.Range("J124").Value = 30000
$ws.Range("L124").Value = 30000
$ws.Range("N124").Value = -39820

$ws.Range("H130").Value = 7429
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 7429
$ws.Range("K130").Value = 0
$ws.Range("L130").ClearContents()
$ws.Range("M130").Value = 7429
$ws.Range("N130").Value = -17469

$ws.Range("H140").Value = 23071
$ws.Range("I140").Value = 19999
$ws.Range("J140").Value = 24299.8
$ws.Range("K140").Value = 19999
$ws.Range("L140").Value = 24299.8
$ws.Range("M140").Value = -14819
$ws.Range("N140").Value = -34659.8

$ws.Range("H141").Value = 101583.336
$ws.Range("I141").Value = 50000
$ws.Range("J141").Value = 127375
$ws.Range("K141").Value = 50000
$ws.Range("L141").Value = 127375
$ws.Range("M141").Value = -44820
$ws.Range("N141").Value = -137735

Write-Output "done"